# Apply updated CircaDB / CircadiPy cosinor analysis results to the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
$ws.Range("G2").Value = 0.0105869880954601
$ws.Range("H2").Value = 0.0344019173187577
$ws.Range("K2").Value = 5.009101870520533
$ws.Range("L2").Value = "[0.6019064606217128, 9.416297280419354]"
$ws.Range("M2").Value = 0.02612383847921929
$ws.Range("N2").Value = 0.02612383847921929
$ws.Range("O2").Value = -0.6918422260157699
$ws.Range("P2").Value = "[-1.471737098979002, 0.08805264694746207]"
$ws.Range("Q2").Value = 0.08176146630394387
$ws.Range("R2").Value = 0.08176146630394387
$ws.Range("S2").Value = 13.16713757656642
$ws.Range("T2").Value = "[10.823959610816654, 15.510315542316178]"
$ws.Range("W2").Value = 2.86176176176183
$ws.Range("X2").Value = -0.364224224224234
$ws.Range("Y2").Value = 6.087747747747894

# --- Row 3 updates ---
$ws.Range("E3").Value = 22.39000000000006
$ws.Range("G3").Value = 0.000110218450380617
$ws.Range("H3").Value = 0.001455160011209773
$ws.Range("K3").Value = 6.494567474846266
$ws.Range("L3").Value = "[2.515576838124659, 10.473558111567874]"
$ws.Range("M3").Value = 0.001483482819091497
$ws.Range("N3").Value = 0.002966965638182995
$ws.Range("O3").Value = 0.7107106503616558
$ws.Range("P3").Value = "[0.1446579199851179, 1.2767633807381937]"
$ws.Range("Q3").Value = 0.01408276315501844
$ws.Range("R3").Value = 0.02816552631003688
$ws.Range("S3").Value = 14.37455453811905
$ws.Range("T3").Value = "[12.282422744892955, 16.46668633134515]"
$ws.Range("W3").Value = 19.85739739739745
$ws.Range("X3").Value = 17.84028028028033
$ws.Range("Y3").Value = 21.87451451451457
